$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (PHP acronym question)
$ws.Range("A2").Value = "What is the acronym of PHP"
$ws.Range("B2").Value = "Hypertext Preprocessor"
$ws.Range("C2").Value = "Preprocessor Hypertext Page"

# Update row 3 (PHP stable version question)
$ws.Range("A3").Value = "Which is the current stable version of php"
$ws.Range("B3").Value = "Php-5.5.9"
$ws.Range("C3").Value = "Php-5.4.5"

# Remove row 4 entirely
$ws.Range("A4:C4").EntireRow.Delete()

# Adjust column widths for columns A and B
# (values pre-compensated for the engine's pixel-rounding of ColumnWidth
#  so the resulting stored width lands as close as possible to the
#  target 30.58988764044944 / 16.28988764044944)
$ws.Columns.Item(1).ColumnWidth = 29.857142857142858
$ws.Columns.Item(2).ColumnWidth = 15.571428571428571
